$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") rows 2 through 469 all currently hold the serial
# date value 45204 (2023-10-05). Update them to 45205 (2023-10-06),
# preserving the existing cell formatting/style.
$lastRow = 469
$range = $ws.Range("C2:C469")

$values = New-Object 'object[,]' $($lastRow - 1), 1
for ($i = 1; $i -le ($lastRow - 1); $i++) {
    $values[$i - 1, 0] = 45205
}

$range.Value = $values
